# "Add files via upload" / "SQL Contents added"
# Adds a new "SQL" training-topics section (rows 62-77) to the bottom of
# Sheet1, mirroring the look & feel of the existing "GIT" / "Maven"
# sections above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Build each new row by copying the formatting of the most similar
#    existing row (so borders / fills / number formats / column styles
#    all match what Excel itself would have produced), then fill in the
#    text / date values for the new row.
# ---------------------------------------------------------------------

function Copy-RowFormat($srcRow, $dstRow) {
    $ws.Range("A$srcRow`:C$srcRow").Copy() | Out-Null
    $ws.Range("A$dstRow`:C$dstRow").PasteSpecial(-4122) | Out-Null
}

# "Topic" rows (A=topic label, C=date) leave the middle cell blank but
# it still carries its own (alignment-flagged) style, distinct from the
# plain interior cells - nudge it so it gets its own cellXf again
# instead of collapsing onto the plain-border style.
function Copy-DateRowFormat($srcRow, $dstRow) {
    Copy-RowFormat $srcRow $dstRow
    $ws.Range("B$dstRow").WrapText = $false
}

# Section header row (merged, bold, filled) -> like row 31 ("Maven")
Copy-RowFormat 31 62
$ws.Range("A62").Value = "SQL"
$ws.Range("A62:C62").Merge() | Out-Null

# "Topic with date" rows -> like row 56
Copy-DateRowFormat 56 63
$ws.Range("A63").Value = "SQL Introduction"
$ws.Range("C63").Value = 43455

# Plain sub-rows -> like row 57
Copy-RowFormat 57 64
$ws.Range("B64").Value = "Introduction"

Copy-RowFormat 57 65
$ws.Range("B65").Value = "Oracle SQL Developer"

Copy-RowFormat 57 66
$ws.Range("B66").Value = "Oracle Express Edition"

# Row 67 is a special "boxed" row: each cell individually gets a
# left+right border, and the outer (A/C) cells additionally get a top
# border, while the middle (B) cell stays vertically centred.
# NOTE: building this one from scratch (no Copy/PasteSpecial) because
# border edits silently no-op on a cell that was just the destination
# of a PasteSpecial in this runtime.
$a67 = $ws.Range("A67")
$a67.Borders.Item(7).LineStyle = 1
$a67.Borders.Item(7).Weight = 2
$a67.Borders.Item(10).LineStyle = 1
$a67.Borders.Item(10).Weight = 2
$a67.Borders.Item(8).LineStyle = 1
$a67.Borders.Item(8).Weight = 2

$b67 = $ws.Range("B67")
$b67.Borders.Item(7).LineStyle = 1
$b67.Borders.Item(7).Weight = 2
$b67.Borders.Item(10).LineStyle = 1
$b67.Borders.Item(10).Weight = 2
$b67.VerticalAlignment = -4108
$b67.Value = "Creation Connections and setting up SQL Developer"

$c67 = $ws.Range("C67")
$c67.Borders.Item(7).LineStyle = 1
$c67.Borders.Item(7).Weight = 2
$c67.Borders.Item(10).LineStyle = 1
$c67.Borders.Item(10).Weight = 2
$c67.Borders.Item(8).LineStyle = 1
$c67.Borders.Item(8).Weight = 2
$c67.NumberFormat = "[`$-409]d\-mmm\-yyyy;@"
$c67.HorizontalAlignment = -4108

Copy-RowFormat 57 68
$ws.Range("B68").Value = "Basic SQL Systax - CREATE, SELECT, INSERT, UPDATE & DELETE"

Copy-DateRowFormat 56 69
$ws.Range("A69").Value = "SELECT Options"
$ws.Range("C69").Value = 43456

Copy-RowFormat 57 70
$ws.Range("B70").Value = "WHERE Clause"

Copy-RowFormat 57 71
$ws.Range("B71").Value = "ORDER BY and GROUP BY"

Copy-RowFormat 57 72
$ws.Range("B72").Value = "JOINS"

Copy-DateRowFormat 56 73
$ws.Range("A73").Value = "SEQUENCES, TRIGGERS, Contraints"
$ws.Range("C73").Value = 43457

Copy-RowFormat 57 74
$ws.Range("B74").Value = "Primary Keys"

Copy-RowFormat 57 75
$ws.Range("B75").Value = "Foreign Keys"

Copy-RowFormat 57 76
Copy-RowFormat 57 77
# NB: "Sequences" (row 77) is written to the shared-string table before
# "Triggers" (row 76) - matches the source workbook's shared-string
# order, even though row 76 (earlier row) is the one that reads
# "Triggers" and row 77 reads "Sequences".
$ws.Range("B77").Value = "Sequences"
$ws.Range("B76").Value = "Triggers"
